$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1469.125
$ws.Range("I98").Value = 1477.8064
$ws.Range("K98").Value = 1477.8064
$ws.Range("M98").Value = 20.19360000000006
$ws.Range("H112").Value = 3214.2856
$ws.Range("J112").Value = 3706.25
$ws.Range("L112").Value = 11118.75
$ws.Range("N112").Value = -13334.75
$ws.Range("H122").Value = 1469.125
$ws.Range("I122").Value = 1477.8064
$ws.Range("K122").Value = 4433.4192
$ws.Range("M122").Value = -1983.4192
$ws.Range("H138").Value = 3200.76
$ws.Range("I138").Value = 1646.2273
$ws.Range("J138").Value = 3639.218
$ws.Range("K138").Value = 4938.6819
$ws.Range("L138").Value = 10917.654
$ws.Range("M138").Value = 201.3181000000004
$ws.Range("N138").Value = -21197.654
$ws.Range("H141").Value = 821.6667
$ws.Range("I141").Value = 770.7143
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 2312.1429
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 2867.8571
$ws.Range("N141").Value = -13360
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22702.855
$ws.Range("I32").Value = 16462.34
$ws.Range("J32").Value = 59365.875
$ws.Range("K32").Value = 16462.34
$ws.Range("L32").Value = 59365.875
$ws.Range("M32").Value = -16175.34
$ws.Range("N32").Value = -59939.875
$ws.Range("H132").Value = 1556.8793
$ws.Range("I132").Value = 1155.7675
$ws.Range("K132").Value = 3467.3025
$ws.Range("M132").Value = -937.3024999999998
$ws.Range("H140").Value = 84321.75
$ws.Range("J140").Value = 84321.75
$ws.Range("L140").Value = 84321.75
$ws.Range("N140").Value = -94681.75
$ws.Range("H141").Value = 52543
$ws.Range("J141").Value = 52543
$ws.Range("L141").Value = 52543
$ws.Range("N141").Value = -62903
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1142.4375
$ws.Range("I16").Value = 991.3570999999999
$ws.Range("J16").Value = 2200
$ws.Range("K16").Value = 991.3570999999999
$ws.Range("L16").Value = 2200
$ws.Range("M16").Value = -704.3570999999999
$ws.Range("N16").Value = -2774
$ws.Range("H31").Value = 2719.1428
$ws.Range("I31").Value = 2276.1875
$ws.Range("J31").Value = 4136.6
$ws.Range("K31").Value = 2276.1875
$ws.Range("L31").Value = 4136.6
$ws.Range("M31").Value = -1981.1875
$ws.Range("N31").Value = -4726.6
$ws.Range("H34").Value = 2719.1428
$ws.Range("I34").Value = 2276.1875
$ws.Range("J34").Value = 4136.6
$ws.Range("K34").Value = 2276.1875
$ws.Range("L34").Value = 4136.6
$ws.Range("M34").Value = -2074.1875
$ws.Range("N34").Value = -4540.6
$ws.Range("H113").Value = 1142.4375
$ws.Range("I113").Value = 991.3570999999999
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 991.3570999999999
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 1178.6429
$ws.Range("N113").Value = -6540
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 9311.111000000001
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9311.111000000001
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 27933.333
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -28789.333
$ws.Range("H91").Value = 9311.111000000001
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9311.111000000001
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 27933.333
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -30897.333
$ws.Range("H107").Value = 519997.2
$ws.Range("I107").Value = 2863.25
$ws.Range("K107").Value = 8589.75
$ws.Range("M107").Value = -6669.75
$ws.Range("H116").Value = 2150
$ws.Range("I116").Value = 1925
$ws.Range("J116").Value = 2420
$ws.Range("K116").Value = 5775
$ws.Range("L116").Value = 7260
$ws.Range("M116").Value = -2333
$ws.Range("N116").Value = -14144
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4279.857
$ws.Range("I70").Value = 3975
$ws.Range("J70").Value = 4401.8
$ws.Range("K70").Value = 3975
$ws.Range("L70").Value = 4401.8
$ws.Range("M70").Value = -3705
$ws.Range("N70").Value = -4941.8
$ws.Range("H73").Value = 4279.857
$ws.Range("I73").Value = 3975
$ws.Range("J73").Value = 4401.8
$ws.Range("K73").Value = 3975
$ws.Range("L73").Value = 4401.8
$ws.Range("M73").Value = -3039
$ws.Range("N73").Value = -6273.8
$ws.Range("H102").Value = 1748.5714
$ws.Range("I102").Value = 1729.2307
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1729.2307
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -107.2307000000001
$ws.Range("N102").Value = -5244
$ws.Range("H126").Value = 3430.4443
$ws.Range("I126").Value = 2483.8572
$ws.Range("J126").Value = 6743.5
$ws.Range("K126").Value = 7451.571599999999
$ws.Range("L126").Value = 20230.5
$ws.Range("M126").Value = -4981.571599999999
$ws.Range("N126").Value = -25170.5
$ws.Range("H132").Value = 3267.4092
$ws.Range("I132").Value = 3248.8
$ws.Range("J132").Value = 3307.2856
$ws.Range("K132").Value = 9746.400000000001
$ws.Range("L132").Value = 9921.856800000001
$ws.Range("M132").Value = -7216.400000000001
$ws.Range("N132").Value = -14981.8568
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2650.1304
$ws.Range("I68").Value = 2558.5
$ws.Range("J68").Value = 2980
$ws.Range("K68").Value = 2558.5
$ws.Range("L68").Value = 2980
$ws.Range("M68").Value = -1809.5
$ws.Range("N68").Value = -4478
$ws.Range("H71").Value = 2650.1304
$ws.Range("I71").Value = 2558.5
$ws.Range("J71").Value = 2980
$ws.Range("K71").Value = 12792.5
$ws.Range("L71").Value = 14900
$ws.Range("M71").Value = -9048.5
$ws.Range("N71").Value = -22388
$ws.Range("H122").Value = 2645.2
$ws.Range("I122").Value = 2388
$ws.Range("J122").Value = 3031
$ws.Range("K122").Value = 7164
$ws.Range("L122").Value = 9093
$ws.Range("M122").Value = -4714
$ws.Range("N122").Value = -13993
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 360.16666
$ws.Range("I113").Value = 287.5
$ws.Range("J113").Value = 396.5
$ws.Range("K113").Value = 862.5
$ws.Range("L113").Value = 1189.5
$ws.Range("M113").Value = 1307.5
$ws.Range("N113").Value = -5529.5
$ws.Range("H126").Value = 844.4666999999999
$ws.Range("I126").Value = 819.0714
$ws.Range("K126").Value = 2457.2142
$ws.Range("M126").Value = 12.78579999999965
